$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @("平潭发展", "平潭发展", "平潭发展"),
    @("福龙马",   "三六零",   "和而泰"),
    @("三六零",   "多氟多",   "海峡创新"),
    @("天际股份", "福龙马",   "福龙马"),
    @("多氟多",   "天际股份", "三六零"),
    @("神州信息", "工业富联", "时空科技"),
    @("和而泰",   "和而泰",   "大众公用"),
    @("海峡创新", "神州信息", "多氟多"),
    @("大众公用", "江特电机", "天际股份"),
    @("工业富联", "利欧股份", "振德医疗"),
    @("中际旭创", "大众公用", "三花智控"),
    @("利欧股份", "海峡创新", "神州信息"),
    @("三花智控", "福石控股", "天赐材料"),
    @("万向钱潮", "万向钱潮", "工业富联"),
    @("江特电机", "舒泰神",   "万向钱潮"),
    @("山子高科", "山子高科", "山子高科"),
    @("舒泰神",   "三花智控", "合富中国"),
    @("天赐材料", "蓝色光标", "舒泰神"),
    @("东方精工", "中际旭创", "利欧股份"),
    @("荣科科技", "天赐材料", "粤传媒")
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    $ws.Cells.Item($row, 1).Value = $rowVals[0]
    $ws.Cells.Item($row, 2).Value = $rowVals[1]
    $ws.Cells.Item($row, 3).Value = $rowVals[2]
}
